$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Ptprc","Mrc1","ECs",2,0.6666666666666666,0.1546356666666667,0.463907,0.0006541814359458435,0.0006541814359458435,3,1,0.316305,0.948915,0.002355313614469985,0.002355313614469984,0.048912034545,0.440208310905,[double]"1.54080244241677e-06",[double]"1.540802442416769e-06"),
    @("ECs","Ptprc","Mrc1","FAPs",2,0.6666666666666666,0.1546356666666667,0.463907,0.0006541814359458435,0.0006541814359458435,3,1,0.2150266666666667,0.64508,0.001601161016974437,0.001601161016974436,0.03325079195111112,0.29925712756,[double]"1.047449813264844e-06",[double]"1.047449813264844e-06"),
    @("ECs","Ptprc","Mrc1","M2",2,0.6666666666666666,0.1546356666666667,0.463907,0.0006541814359458435,0.0006541814359458435,3,1,130.539482,391.618446,0.9720409705204137,0.9720409705204136,20.18605982539133,181.674538428522,0.0006358911578932356,0.0006358911578932356),
    @("ECs","Ptprc","Mrc1","sCs",2,0.6666666666666666,0.1546356666666667,0.463907,0.0006541814359458435,0.0006541814359458435,3,1,3.223404333333333,9.670213,0.02400255484814197,0.02400255484814196,0.4984532780212222,4.486079502191,[double]"1.570202579692638e-05",[double]"1.570202579692638e-05"),
    @("FAPs","Ptprc","Mrc1","ECs",3,1,0.3253436666666666,0.976031,0.001376356384162467,0.001376356384162467,3,1,0.316305,0.948915,0.002355313614469985,0.002355313614469984,0.102907828485,0.926170456365,[double]"3.241750929980539e-06",[double]"3.241750929980538e-06"),
    @("FAPs","Ptprc","Mrc1","FAPs",3,1,0.3253436666666666,0.976031,0.001376356384162467,0.001376356384162467,3,1,0.2150266666666667,0.64508,0.001601161016974437,0.001601161016974436,0.06995756416444444,0.62961807748,[double]"2.203768187784834e-06",[double]"2.203768187784834e-06"),
    @("FAPs","Ptprc","Mrc1","M2",3,1,0.3253436666666666,0.976031,0.001376356384162467,0.001376356384162467,3,1,130.539482,391.618446,0.9720409705204137,0.9720409705204136,42.47019371864733,382.231743467826,0.001337874795443252,0.001337874795443252),
    @("FAPs","Ptprc","Mrc1","sCs",3,1,0.3253436666666666,0.976031,0.001376356384162467,0.001376356384162467,3,1,3.223404333333333,9.670213,0.02400255484814197,0.02400255484814196,1.048714184955889,9.438427664603001,[double]"3.303606960144997e-05",[double]"3.303606960144996e-05"),
    @("M2","Ptprc","Mrc1","ECs",3,1,231.9004236666667,695.701271,0.9810476161216116,0.9810476161216116,3,1,0.316305,0.948915,0.002355313614469985,0.002355313614469984,73.35126350788501,660.161371570965,0.002310674806694555,0.002310674806694554),
    @("M2","Ptprc","Mrc1","FAPs",3,1,231.9004236666667,695.701271,0.9810476161216116,0.9810476161216116,3,1,0.2150266666666667,0.64508,0.001601161016974437,0.001601161016974436,49.86477509963112,448.78297589668,0.001570815198729627,0.001570815198729626),
    @("M2","Ptprc","Mrc1","M2",3,1,231.9004236666667,695.701271,0.9810476161216116,0.9810476161216116,3,1,130.539482,391.618446,0.9720409705204137,0.9720409705204136,30272.16118102721,272449.4506292448,0.9536184769015896,0.9536184769015895),
    @("M2","Ptprc","Mrc1","sCs",3,1,231.9004236666667,695.701271,0.9810476161216116,0.9810476161216116,3,1,3.223404333333333,9.670213,0.02400255484814197,0.02400255484814196,747.5088305489693,6727.579474940723,0.02354764921459791,0.0235476492145979),
    @("sCs","Ptprc","Mrc1","ECs",3,1,3.999992666666667,11.999978,0.01692184605828007,0.01692184605828007,3,1,0.316305,0.948915,0.002355313614469985,0.002355313614469984,1.26521768043,11.38695912387,[double]"3.98562544030323e-05",[double]"3.985625440303228e-05"),
    @("sCs","Ptprc","Mrc1","FAPs",3,1,3.999992666666667,11.999978,0.01692184605828007,0.01692184605828007,3,1,0.2150266666666667,0.64508,0.001601161016974437,0.001601161016974436,0.8601050898044446,7.74094580824,[double]"2.709460024376058e-05",[double]"2.709460024376057e-05"),
    @("sCs","Ptprc","Mrc1","M2",3,1,3.999992666666667,11.999978,0.01692184605828007,0.01692184605828007,3,1,130.539482,391.618446,0.9720409705204137,0.9720409705204136,522.1569707104653,4699.412736394188,0.0164487276654876,0.0164487276654876),
    @("sCs","Ptprc","Mrc1","sCs",3,1,3.999992666666667,11.999978,0.01692184605828007,0.01692184605828007,3,1,3.223404333333333,9.670213,0.02400255484814197,0.02400255484814196,12.89359369503489,116.042343255314,0.0004061675381456823,0.0004061675381456822)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}
